$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1: copy formatting from B1 (bold/border/centered header style),
# then set its text so it reuses the existing header style record.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Coord: normal vector scan"

# Updated column B values (recomputed angles) and new column C (scan normal vector text)
$ws.Range("B2").Value = 0.3560034938780682
$ws.Range("C2").Value = "[0.         0.31908379 0.94772651]"
$ws.Range("B3").Value = 1.516571111967983
$ws.Range("C3").Value = "[-0.43876148  0.52155555  0.73175691]"
$ws.Range("B4").Value = 0.1963039095234749
$ws.Range("C4").Value = "[-0.00119071  0.01282932  0.99991699]"
$ws.Range("B5").Value = 0.6034302761849968
$ws.Range("C5").Value = "[ 2.54597526e-04  2.84382834e-01 -9.58710769e-01]"
$ws.Range("B6").Value = 1.694921828085868
$ws.Range("C6").Value = "[0.72390296 0.28901965 0.62644405]"
$ws.Range("B7").Value = 0.982589045564234
$ws.Range("C7").Value = "[-0.73129183 -0.27139518  0.62574509]"
$ws.Range("B8").Value = 0.588816101442688
$ws.Range("C8").Value = "[0.         0.31523022 0.94901523]"
$ws.Range("B9").Value = 0.8121553091191543
$ws.Range("C9").Value = "[ 0.         -0.31152858  0.95023678]"
$ws.Range("B10").Value = 2.555862846884493
$ws.Range("C10").Value = "[-0.72303144  0.24963508  0.64413342]"
$ws.Range("B11").Value = 1.827680287833167
$ws.Range("C11").Value = "[ 0.72947005 -0.27563522  0.62601811]"
$ws.Range("B12").Value = 0.6637576211436756
$ws.Range("C12").Value = "[-2.16747601e-04 -2.83372235e-01 -9.59009974e-01]"
$ws.Range("B13").Value = 1.457736114015883
$ws.Range("C13").Value = "[-0.72462879 -0.29607244  0.62229754]"
$ws.Range("B14").Value = 2.603716598398236
$ws.Range("C14").Value = "[0.71326391 0.26845507 0.64744611]"
$ws.Range("B15").Value = 1.18696440254348
$ws.Range("C15").Value = "[ 0.         -0.30530584  0.95225435]"
